# Add a new "Unindexed" sheet (solver run metadata) after the existing
# sheets, matching the commit's refactor of profit computation / price
# extraction which now also emits solver diagnostics.

$wb = $excel.ActiveWorkbook

$firstSheet = $wb.Worksheets.Item(1)
$lastSheet  = $wb.Worksheets.Item($wb.Worksheets.Count)

# New sheet goes after the last existing sheet (sPaths).
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Unindexed"

# Header row.
$ws.Range("A1").Value = "solver_status"
$ws.Range("B1").Value = "termination_condition"
$ws.Range("C1").Value = "objective_value"
$ws.Range("D1").Value = "final_gap"
$ws.Range("E1").Value = "execution_time"
$ws.Range("F1").Value = "lower_bound"
$ws.Range("G1").Value = "upper_bound"

# Reuse the workbook's existing bold/bordered/centered header style (the
# same one already used for the header rows on sIntersections / sPaths)
# instead of building a fresh style from scratch.
$firstSheet.Range("A1").Copy() | Out-Null
$ws.Range("A1:G1").PasteSpecial(-4122) | Out-Null

# Data row.
$ws.Range("A2").Value = "aborted"
$ws.Range("B2").Value = "maxTimeLimit"
$ws.Range("C2").Value = 9.843673651466128
$ws.Range("D2").Value = 0.9999999999999994
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = 0.000000000000005329070518200751
$ws.Range("G2").Value = 9.843673651466128

# Restore the original active sheet/selection so the only visible change
# is the newly added sheet.
$firstSheet.Activate() | Out-Null
$firstSheet.Range("A1").Select() | Out-Null
